# Weekly data refresh: insert two new rows of Kiwi price data at the top
# of the data block (row 38), pushing the existing rows down by two and
# growing the used range from A1:T140 to A1:T142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 38; this shifts rows 38:140 down to 40:142,
# carrying their formatting (and the dimension) along automatically.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()

function Set-KiwiRow($row, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $precioKg) {
    $ws.Cells.Item($row, 1).Value = 7
    $ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($row, 3).Value = "Ñuble"
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = 16
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100101007
    $ws.Cells.Item($row, 10).Value = "Kiwi"
    $ws.Cells.Item($row, 11).Value = "Hayward"
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $pmin
    $ws.Cells.Item($row, 15).Value = $pmax
    $ws.Cells.Item($row, 16).Value = $pprom
    $ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = 18
}

# New row 38 ("Primera" quality, week of 2021-09-14)
Set-KiwiRow 38 44453 "Primera" 100 12000 12500 12250 681

# New row 39 ("Segunda" quality, week of 2021-09-14)
Set-KiwiRow 39 44453 "Segunda" 60 11000 11500 11250 625
